$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.110.82"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.835.17"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'414.55"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "'133.05"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "3.824.02"
$ws.Range("E7").Value = "  +3.63%  "
$ws.Range("D8").Value = "'0.620"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.747"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "'0.173"
$ws.Range("E11").Value = "  -3.95%  "
$ws.Range("D12").Value = "'0.0000380"
$ws.Range("E12").Value = "  -2.63%  "
$ws.Range("D13").Value = "'41.35"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("D14").Value = "4.468.37"
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("D15").Value = "'10.08"
$ws.Range("E15").Value = "  -5.09%  "
$ws.Range("D16").Value = "'14.93"
$ws.Range("E16").Value = "  +14.46%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.138"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.803.82"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'19.62"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").Value = "67.492.31"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "'417.81"
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "'15.03"
$ws.Range("E23").Value = "  -8.67%  "
$ws.Range("D24").Value = "'86.67"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").Value = "'3.10"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'36.71"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "'5.74"
$ws.Range("E27").Value = "  +14.58%  "
$ws.Range("D28").Value = "'3.17"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("D29").Value = "'9.60"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("D30").Value = "'691.76"
$ws.Range("E30").Value = "  +5.86%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'12.55"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.122"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "'2.74"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "'7.28"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "'0.154"
$ws.Range("E35").Value = "  -7.42%  "
$ws.Range("D36").Value = "'39.33"
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'55.58"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "0.0₃0786"
$ws.Range("E39").Value = "  +7.67%  "
$ws.Range("D40").Value = "'0.0463"
$ws.Range("E40").Value = "  -6.07%  "
$ws.Range("D41").Value = "'3.09"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  -8.17%  "
$ws.Range("D44").Value = "'27.46"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'148.27"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "'3.34"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'4.45"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.17"
$ws.Range("E48").Value = "  +18.91%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'2.11"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'2.87"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'2.61"
$ws.Range("E51").Value = "  +0.40%  "
